# results_crossentropy.xlsx - "unified report and thesis documents"
#
# The "Type - Image augmentation" block (rows 15-20) was missing the raw
# per-run Dev/Test accuracy figures for Run#2..Run#5 (rows 17-20, columns
# B/C) that every other block on the sheet already has. This fills those
# in, which lets the existing MIN/MAX/AVERAGE formulas in D16:I16 (and the
# derived spread formulas in F17/I17) recalculate to their correct values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Run#2 (row 17)
$ws.Range("B17").Value = 0.75482912332838004
$ws.Range("C17").Value = 0.764589515331355

# Run#3 (row 18)
$ws.Range("B18").Value = 0.75928677563150004
$ws.Range("C18").Value = 0.74085064292779401

# Run#4 (row 19)
$ws.Range("B19").Value = 0.75928677563150004
$ws.Range("C19").Value = 0.75865479723046403

# Run#5 (row 20)
$ws.Range("B20").Value = 0.76374442793462105
$ws.Range("C20").Value = 0.73590504451038496

# Recalculate so the MIN/MAX/AVERAGE/spread formulas pick up the new data.
$excel.CalculateFull()

# Match the author's final selection state (cell F16).
$ws.Range("F16").Select()
